$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSGP")

# Row 16: Gross Margin
$ws.Range("D16").Value = 0.7795
$ws.Range("E16").Value = 0.7853
$ws.Range("F16").Value = 0.7892
$ws.Range("G16").Value = 0.7934

# Row 20: Free Cash Flow Margin
$ws.Range("D20").Value = 0.277
$ws.Range("E20").Value = 0.2917
$ws.Range("F20").Value = 0.276
$ws.Range("G20").Value = 0.294

# Row 28: EBITDA Margin
$ws.Range("D28").Value = 0.2995
$ws.Range("E28").Value = 0.3219
$ws.Range("F28").Value = 0.3278
$ws.Range("G28").Value = 0.3565

# Row 29: Operating Cash Flow Margin
$ws.Range("D29").Value = 0.3047
$ws.Range("E29").Value = 0.3211
$ws.Range("F29").Value = 0.306
$ws.Range("G29").Value = 0.3271
